$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the generated Date value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(8, 2).Value = "2026-01-23T08:28:04+00:00"

# --- Mapping Table 0: split "dateDebutFin"/"effectiveTime" row into two rows
#     (dateDebut -> effectiveTime.low, dateFin -> effectiveTime.high) ---
$map0 = $wb.Worksheets.Item("Mapping Table 0")

# Insert a fresh row right after row 7, pushing the remaining rows down.
$map0.Rows.Item(8).Insert()

# Give the new row the same look (borders/wrap) as the surrounding data rows.
$map0.Range("A7:E7").Copy()
$map0.Range("A8:E8").PasteSpecial(-4122)

# Row 7 now represents the "start date" half of the mapping.
$map0.Cells.Item(7, 1).Value = "FRLMEffetIndesirable.dateDebut"
$map0.Cells.Item(7, 3).Value = "equivalent"
$map0.Cells.Item(7, 4).Value = "FRCDAEffetIndesirable.effectiveTime.low"

# Row 8 (new) represents the "end date" half of the mapping.
$map0.Cells.Item(8, 1).Value = "FRLMEffetIndesirable.dateFin"
$map0.Cells.Item(8, 3).Value = "equivalent"
$map0.Cells.Item(8, 4).Value = "FRCDAEffetIndesirable.effectiveTime.high"
